# Auto-generated edit script: refresh crypto price/volume data
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "21.381.57"
$ws.Cells.Item(2, 5).Value = "  -2.90%  "
$ws.Cells.Item(3, 4).Value = "1.529.30"
$ws.Cells.Item(3, 5).Value = "  -1.57%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.005"
$ws.Cells.Item(4, 5).Value = "  +0.20%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "1.003"
$ws.Cells.Item(5, 5).Value = "  +0.06%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "288.11"
$ws.Cells.Item(6, 5).Value = "  -0.96%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3902"
$ws.Cells.Item(7, 5).Value = "  -0.52%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3171"
$ws.Cells.Item(8, 5).Value = "  -1.32%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "42.59"
$ws.Cells.Item(9, 5).Value = "  -2.15%  "
$ws.Cells.Item(10, 5).Value = "  -1.72%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.067"
$ws.Cells.Item(11, 5).Value = "  -0.57%  "
$ws.Cells.Item(12, 5).Value = "  +0.04%  "
$ws.Cells.Item(13, 5).Value = "  +1.13%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "18.15"
$ws.Cells.Item(14, 5).Value = "  -2.89%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "6.477"
$ws.Cells.Item(15, 5).Value = "  -2.01%  "
$ws.Cells.Item(16, 4).Value = "1.529.45"
$ws.Cells.Item(16, 5).Value = "  -1.44%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001084"
$ws.Cells.Item(17, 5).Value = "  -3.71%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.06599"
$ws.Cells.Item(18, 5).Value = "  +0.13%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "82.86"
$ws.Cells.Item(19, 5).Value = "  -0.70%  "
$ws.Cells.Item(20, 5).Value = "  +0.10%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.089"
$ws.Cells.Item(21, 5).Value = "  -2.99%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "15.37"
$ws.Cells.Item(22, 5).Value = "  -0.93%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.83"
$ws.Cells.Item(23, 5).Value = "  -3.68%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.372"
$ws.Cells.Item(24, 5).Value = "  +0.32%  "
$ws.Cells.Item(25, 4).Value = "21.446.71"
$ws.Cells.Item(25, 5).Value = "  -2.67%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.366"
$ws.Cells.Item(26, 5).Value = "  -2.12%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "148.92"
$ws.Cells.Item(27, 5).Value = "  +0.15%  "
$ws.Cells.Item(28, 5).Value = "  -1.23%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "4.822"
$ws.Cells.Item(29, 5).Value = "  -1.41%  "
$ws.Cells.Item(30, 4).Value = "1.702.49"
$ws.Cells.Item(30, 5).Value = "  -1.40%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "116.45"
$ws.Cells.Item(31, 5).Value = "  -2.00%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.022"
$ws.Cells.Item(32, 5).Value = "  +4.61%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.9481"
$ws.Cells.Item(33, 5).Value = "  -4.95%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08000"
$ws.Cells.Item(34, 5).Value = "  -3.64%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "8.449"
$ws.Cells.Item(35, 5).Value = "  -5.96%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.134"
$ws.Cells.Item(36, 5).Value = "  +0.90%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.492"
$ws.Cells.Item(37, 5).Value = "  -7.99%  "
$ws.Cells.Item(38, 2).Value = "Aptos"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "11.24"
$ws.Cells.Item(38, 5).Value = "  +5.45%  "
$ws.Cells.Item(39, 2).Value = "Hedera"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.05906"
$ws.Cells.Item(39, 5).Value = "  -2.80%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.02181"
$ws.Cells.Item(40, 5).Value = "  -3.50%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.2015"
$ws.Cells.Item(41, 5).Value = "  -1.14%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.175"
$ws.Cells.Item(42, 5).Value = "  -2.94%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.5729"
$ws.Cells.Item(44, 5).Value = "  -1.35%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "13.05"
$ws.Cells.Item(45, 5).Value = "  -0.47%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.708"
$ws.Cells.Item(46, 5).Value = "  -1.19%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5539"
$ws.Cells.Item(47, 5).Value = "  -0.48%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.891"
$ws.Cells.Item(48, 5).Value = "  -0.19%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.156"
$ws.Cells.Item(49, 5).Value = "  +2.14%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "115.53"
$ws.Cells.Item(50, 5).Value = "  -2.33%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.06655"
$ws.Cells.Item(51, 5).Value = "  -2.40%  "
